$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A200").Value = "Hello"
